$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# ---------------------------------------------------------------------------
# 1) All existing "Results" values (E2:E19) flip from PASS to SKIP
# ---------------------------------------------------------------------------
$ws.Range("E2:E19").Value = "SKIP"

# ---------------------------------------------------------------------------
# 2) Append three brand-new test cases (rows 20-22), copying the formatting
#    (borders / fill / wrap) from the last existing data row (19) first so
#    the new rows look consistent with the rest of the table.
# ---------------------------------------------------------------------------
$ws.Range("A19:E19").Copy()
$ws.Range("A20:E22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 20 - TestCase_E19
$ws.Range("A20").Value = "TestCase_E19"
$ws.Range("B20").Value = "OPQA-288"
$ws.Range("C20").Value = "Verify that following fields are getting displayed for each article in the watchlist page:`na)Times cited`nb)Comments"
$ws.Range("D20").Value = "Y"
$ws.Range("E20").Value = "SKIP"

# Row 21 - TestCase_E20
$ws.Range("A21").Value = "TestCase_E20"
$ws.Range("B21").Value = "OPQA-290"
$ws.Range("C21").Value = "Verify that following fields are getting displayed for each article in the watchlist page:`na)Times cited`nb)Comments"
$ws.Range("D21").Value = "Y"
$ws.Range("E21").Value = "SKIP"

# Row 22 - TestCase_E21
$ws.Range("A22").Value = "TestCase_E21"
$ws.Range("B22").Value = "OPQA-291"
$ws.Range("C22").Value = "Verify that following fields are getting displayed for each post in the watchlist page:`na)Likes`nb)Comments"
$ws.Range("D22").Value = "Y"
$ws.Range("E22").Value = "PASS"

# Match the row height used by the author for the new, taller (wrapped) rows
$ws.Range("A20:E22").RowHeight = 45

# ---------------------------------------------------------------------------
# 3) Update the sheet's selection/view to reflect where the author left off
# ---------------------------------------------------------------------------
[void]$ws.Range("D2:D22").Select()

Write-Host "done"
